$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates derived from the authoritative diff (cryptos list refresh).
# A handful of the new price strings are unambiguous numerals (e.g. "1.00",
# "604.07"); Excel would silently coerce those to numbers on a plain .Value
# assignment, so they get a leading apostrophe to force literal text, matching
# every other cell in this column (all originally stored as text).
$ws.Range("D2").Value = "66.187.50"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "3.552.31"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'604.07"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'146.55"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("D7").Value = "3.552.57"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").Value = "'0.496"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").Value = "'0.133"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("D11").Value = "'7.80"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "'0.408"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "4.155.55"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "'0.0000201"
$ws.Range("E14").Value = "  -2.77%  "
$ws.Range("D15").Value = "'29.06"
$ws.Range("E15").Value = "  -3.94%  "
$ws.Range("D16").Value = "3.555.91"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("D18").Value = "66.197.26"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'11.05"
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("D20").Value = "'6.22"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").Value = "'14.71"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").Value = "'418.21"
$ws.Range("E22").Value = "  -3.07%  "
$ws.Range("D23").Value = "'0.603"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("D24").Value = "'77.95"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("D25").Value = "3.694.39"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'0.0000117"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").Value = "'9.12"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.47"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'7.88"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "3.550.29"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").Value = "'0.155"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").Value = "'24.52"
$ws.Range("E34").Value = "  -3.66%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'7.61"
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("D37").Value = "'1.31"
$ws.Range("E37").Value = "  -9.50%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.61"
$ws.Range("E38").Value = "  -6.92%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'5.29"
$ws.Range("E39").Value = "  -5.64%  "
$ws.Range("D40").Value = "'173.67"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").Value = "'0.0823"
$ws.Range("E41").Value = "  -2.85%  "
$ws.Range("D42").Value = "'5.10"
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("D43").Value = "'0.864"
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("D44").Value = "'45.61"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").Value = "'1.81"
$ws.Range("E45").Value = "  -5.84%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'2.44"
$ws.Range("E47").Value = "  -3.81%  "
$ws.Range("D48").Value = "'7.09"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("D49").Value = "'22.75"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("E50").Value = "  -7.12%  "
$ws.Range("D51").Value = "'23.27"
$ws.Range("E51").Value = "  -7.88%  "
